# feat: get taxiing time for each gate
#
# - Rename sheet "mintaxitime " (trailing space) -> "mintaxitime"
# - On the query table "表_表1" (sheet1): hide the per-column AutoFilter
#   drop-down buttons (adds <filterColumn .../> entries) and clear the
#   named table style so the table no longer carries "TableStyleMedium7"
# - Restore the cell selections on each sheet to match the final saved
#   state, leaving "mintaxitime" as the active tab

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("sheet1")
$mintaxi = $wb.Worksheets.Item("mintaxitime ")

# --- table formatting -------------------------------------------------
$lo = $sheet1.ListObjects.Item(1)
$lo.ShowAutoFilterDropDown = $false
$lo.TableStyle = ""

# --- sheet rename -------------------------------------------------------
$mintaxi.Name = "mintaxitime"

# --- view / selection state ---------------------------------------------
# Visit sheet1 first and leave its selection at E6 ...
$sheet1.Activate()
$sheet1.Range("E6").Select()

# ... then return to the mintaxitime sheet so it ends up the active tab,
# with its own selection moved to K20.
$mintaxi.Activate()
$mintaxi.Range("K20").Select()
